$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "51.787.24"
$ws.Range("E2").Value = "  -0.13%  "
$ws.Range("D3").Value = "2.773.10"
$ws.Range("E3").Value = "  -1.70%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "357.83"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -0.61%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "109.31"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -3.65%  "
$ws.Range("E7").Value = "  +2.42%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "1.00"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  +0.08%  "
$ws.Range("E9").Value = "  -1.84%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "39.95"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -4.23%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0849"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -0.11%  "
$ws.Range("E12").Value = "  +0.42%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "19.48"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  -1.90%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.62"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -2.14%  "
$ws.Range("D15").Value = "3.208.79"
$ws.Range("E15").Value = "  -1.83%  "
$ws.Range("D16").Value = "2.767.89"
$ws.Range("E16").Value = "  -2.69%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.915"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  +2.74%  "
$ws.Range("D18").Value = "51.679.48"
$ws.Range("E18").Value = "  -0.40%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.39"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +1.74%  "
$ws.Range("E20").Value = "  -1.45%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.10"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -3.35%  "
$ws.Range("D22").Value = "0.0₃0978"
$ws.Range("E22").Value = "  -0.17%  "
$ws.Range("E23").Value = "  +2.24%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "69.52"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -0.27%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.76"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -1.63%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "26.44"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -1.25%  "
$ws.Range("E27").Value = "  +0.12%  "
$ws.Range("E28").Value = "  -0.85%  "
$ws.Range("E29").Value = "  -0.78%  "
$ws.Range("E30").Value = "  +1.26%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.0473"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +7.06%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "51.22"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +0.90%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "33.92"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +0.67%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.73"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -0.77%  "
$ws.Range("E35").Value = "  +8.52%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0837"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +1.35%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.999"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -0.19%  "
$ws.Range("E38").Value = "  -1.61%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.00"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -4.64%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "17.99"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -1.79%  "
$ws.Range("B41").Value = "Stacks"
$ws.Range("C41").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.59"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +2.27%  "
$ws.Range("B42").Value = "Stellar"
$ws.Range("C42").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.115"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -0.69%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "125.26"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -1.91%  "
$ws.Range("E44").Value = "  -0.89%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "21.99"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -6.87%  "
$ws.Range("D46").Value = "2.055.37"
$ws.Range("E46").Value = "  -0.02%  "
$ws.Range("E47").Value = "  +0.64%  "
$ws.Range("E48").Value = "  -3.67%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "5.70"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  +1.99%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.932"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +0.17%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "8.97"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +0.37%  "
